{"js": "// Replace the date line and the 25 division expressions in the table,\n// per the captured diff (2024-12-03 Tuesday -> 2024-12-04 Wednesday, plus\n// each \"a\u00f7b=\" cell getting a new pair of operands).\nconst replacements = [\n  [\"2024-12-03 Tuesday\", \"2024-12-04 Wednesday\"],\n  [\"92\u00f76=\", \"76\u00f79=\"],\n  [\"83\u00f76=\", \"34\u00f72=\"],\n  [\"40\u00f74=\", \"54\u00f73=\"],\n  [\"76\u00f76=\", \"53\u00f78=\"],\n  [\"58\u00f77=\", \"39\u00f74=\"],\n  [\"22\u00f74=\", \"78\u00f79=\"],\n  [\"77\u00f76=\", \"85\u00f76=\"],\n  [\"85\u00f74=\", \"42\u00f77=\"],\n  [\"99\u00f72=\", \"66\u00f74=\"],\n  [\"90\u00f76=\", \"39\u00f74=\"],\n  [\"48\u00f74=\", \"59\u00f73=\"],\n  [\"97\u00f72=\", \"12\u00f74=\"],\n  [\"74\u00f72=\", \"30\u00f76=\"],\n  [\"75\u00f76=\", \"50\u00f74=\"],\n  [\"20\u00f73=\", \"24\u00f78=\"],\n  [\"97\u00f76=\", \"65\u00f73=\"],\n  [\"22\u00f73=\", \"64\u00f79=\"],\n  [\"35\u00f75=\", \"10\u00f79=\"],\n  [\"98\u00f75=\", \"77\u00f75=\"],\n  [\"99\u00f75=\", \"40\u00f74=\"],\n  [\"82\u00f73=\", \"21\u00f75=\"],\n  [\"23\u00f79=\", \"69\u00f79=\"],\n  [\"39\u00f72=\", \"95\u00f73=\"],\n  [\"30\u00f77=\", \"18\u00f75=\"],\n  [\"91\u00f76=\", \"52\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 division expressions in the table,\n# per the captured diff (2024-12-03 Tuesday -> 2024-12-04 Wednesday, plus\n# each \"a\u00f7b=\" cell getting a new pair of operands).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-12-03 Tuesday\", \"2024-12-04 Wednesday\"),\n    @(\"92\u00f76=\", \"76\u00f79=\"),\n    @(\"83\u00f76=\", \"34\u00f72=\"),\n    @(\"40\u00f74=\", \"54\u00f73=\"),\n    @(\"76\u00f76=\", \"53\u00f78=\"),\n    @(\"58\u00f77=\", \"39\u00f74=\"),\n    @(\"22\u00f74=\", \"78\u00f79=\"),\n    @(\"77\u00f76=\", \"85\u00f76=\"),\n    @(\"85\u00f74=\", \"42\u00f77=\"),\n    @(\"99\u00f72=\", \"66\u00f74=\"),\n    @(\"90\u00f76=\", \"39\u00f74=\"),\n    @(\"48\u00f74=\", \"59\u00f73=\"),\n    @(\"97\u00f72=\", \"12\u00f74=\"),\n    @(\"74\u00f72=\", \"30\u00f76=\"),\n    @(\"75\u00f76=\", \"50\u00f74=\"),\n    @(\"20\u00f73=\", \"24\u00f78=\"),\n    @(\"97\u00f76=\", \"65\u00f73=\"),\n    @(\"22\u00f73=\", \"64\u00f79=\"),\n    @(\"35\u00f75=\", \"10\u00f79=\"),\n    @(\"98\u00f75=\", \"77\u00f75=\"),\n    @(\"99\u00f75=\", \"40\u00f74=\"),\n    @(\"82\u00f73=\", \"21\u00f75=\"),\n    @(\"23\u00f79=\", \"69\u00f79=\"),\n    @(\"39\u00f72=\", \"95\u00f73=\"),\n    @(\"30\u00f77=\", \"18\u00f75=\"),\n    @(\"91\u00f76=\", \"52\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
